$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from F1 to G1, then set header text
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = "spareCosts"

$ws.Range("G2").Value = "{'BOLT': 1.13, 'ELBOW': 533.52, 'PLATE': 153.27, 'SEAL': 824.85, 'SEALANT (310 ML)': 5.8}"
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = "{'DISCHRGR': 357.27}"
$ws.Range("G5").Value = "{'LAMP': 1.35, 'LAMP-GE S-8 28V 643W SC BAYONET': 2.52}"
$ws.Range("G6").Value = "{'ELBOW': 177.84, 'SEAL': 487.89}"
$ws.Range("G7").Value = "{'CORROSION INHIBITING COMPOUND-SPRAY': 39.66, 'FOAM TAPE': 158.13, 'Low Density Sealant': 148.01, 'NON AQUEOUS CLEANER-GENERAL': 13.75, 'SEALANT (130 ML)': 101.44, 'corrosion inhibiting compound': 45.8}"
$ws.Range("G8").Value = "{'ADHESIVE HYSOL EA9309-3 1 QUART PER KIT BMSÂ 5-109 Type I, Class 1': 428.94, 'ALEXIT-FST FILLER 495-14 incl. Hardener 491-14': 803.26, 'FOAM TAPE': 368.96, 'GILLFAB 0.496""X60""X144""': 10205.33, 'INSERT ASSY SLEEVE AND PLUG': 669.24, 'ISOPROPYL ALCOHOL': 7.49, 'Low Density Sealant': 296.01, 'corrosion inhibiting compound': 45.8}"
$ws.Range("G9").Value = "{'SCREW': 20.95, 'STUD': 40.15, 'WASHER': 21.880000000000003}"
$ws.Range("G10").Value = "{'SCREW': 28.43, 'SECTION': 104.01, 'STUD': 400.41, 'WASHER': 6.84}"
$ws.Range("G11").Value = "{'STUD': 9.87, 'WASHER': 2.01}"
$ws.Range("G12").Value = "{'ACCESS DOOR SEALANT': 160.02, 'VALVE': 62.06}"
$ws.Range("G13").Value = "{'FILTER': 153.15}"
$ws.Range("G14").Value = "{'CLAMP': 159.12, 'HOSE': 288.99}"
$ws.Range("G15").Value = "{'BUSHING': 3.58}"
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = "{'CABLE': 224.64, 'LEAD': 224.26999999999998}"
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = "{'SCREW': 1.73, 'SEAL': 215.28}"
$ws.Range("G25").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = "{'SEAL': 215.28}"
$ws.Range("G29").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("G31").Value = "{'ALODINE CHROMATING POWDER': 9.43, 'Demineralized Water.': 1.42, 'Fuel Tank Sealant': 53.69, 'NON AQUEOUS CLEANER-GENERAL': 6.93}"
$ws.Range("G32").Value = 0
$ws.Range("G33").Value = "{'ADHESIVE HYSOL EA9309-3 1 QUART PER KIT BMSÂ 5-109 Type I, Class 1': 428.22, 'Disposable Protective Coverall Safety Work Wear 3M 4540 Medium Body Suit MAKE(KEELGUARD)': 37.69, 'FOAM TAPE': 527.09, 'GILLFAB 0.496""X60""X144""': 10477.47, 'INSERT ASSY SLEEVE AND PLUG': 823.68, 'ISOPROPYL ALCOHOL': 18.69, 'Low Density Sealant': 222.01, 'SEALANT (130 ML)': 114.62, 'corrosion inhibiting compound': 45.78, 'Â Edge Sealing Compound': 1752.55}"
$ws.Range("G34").Value = "{'BLIND RIVET': 4.59, 'FOAM TAPE': 210.83, 'Low Density Sealant': 296.01, 'NUT': 3.7, 'NUT CLIP': 116.06, 'PIN': 0.12000000000000001, 'corrosion inhibiting compound': 45.78}"
$ws.Range("G35").Value = "{'SCREW': 13.46, 'SECTION': 183.35}"
$ws.Range("G36").Value = "{'SCREW': 11.23, 'SECTION': 352.55, 'STUD': 534.02, 'WASHER': 13.07}"
$ws.Range("G37").Value = "{'Access Door Sealant': 88.48, 'VALVE': 82.74}"
$ws.Range("G38").Value = "{'CONTACT': 26.68}"
$ws.Range("G39").Value = "{'DISCHRGR': 608.0}"
$ws.Range("G40").Value = "{'BELLOWS': 596.7}"
$ws.Range("G41").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("G43").Value = "{'LAMP': 33.08}"
$ws.Range("G44").Value = "{'RECEPTCL': 7.38, 'SCREW': 13.0, 'STUD': 18.4, 'WASHER': 4.77}"
$ws.Range("G45").Value = "{'CABLE': 1095.12, 'Fuel Tank Sealant': 53.69, 'LEAD': 232.94}"
$ws.Range("G46").Value = "{'FILTER': 14.65}"
$ws.Range("G47").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("G49").Value = "{'BOLT': 0.19, 'COMPASS LIGHT': 99.98, 'LAMP': 37.85, 'MOULD RELEASE AGENT': 31.82, 'SCREW': 2.15, 'SEALANT': 159.12}"
$ws.Range("G50").Value = 0
$ws.Range("G51").Value = "{'ADHESIVE HYSOL EA9309-3 1 QUART PER KIT BMSÂ 5-109 Type I, Class 1': 424.37, 'ADHESIVE TAPE': 378.14, 'Disposable Protective Coverall Safety Work Wear 3M 4520 Body Suit': 36.0, 'Disposable Protective Coverall Safety Work Wear 3M 4540 Medium Body Suit MAKE(KEELGUARD)': 13.55, 'Fuel Tank Sealant': 256.21, 'GILLFAB 0.496""X60""X144""': 3596.78, 'NUT CLIP': 11.61, 'SCREW': 4.49, 'SECTION': 85.41, 'STUD': 305.59, 'WASHER': 6.42, 'Â Edge Sealing Compound': 553.37}"
$ws.Range("G52").Value = "{'BRAID': 68.33000000000001, 'CABLE': 375.34000000000003, 'LEAD': 189.51, 'SEALANT': 65.23}"
$ws.Range("G53").Value = "{'ADHESIVE TAPE': 189.07, 'CORROSION INHIBITING COMPOUND': 45.49, 'Fuel Tank Sealant': 768.62, 'INSERT ASSY SLEEVE AND PLUG': 134.55, 'METHYL ETHYL KETONE': 6.98, 'NON AQUEOUS CLEANER-GENERAL': 6.88, 'SCREW': 19.14, 'STUD': 557.3499999999999}"
$ws.Range("G54").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("G57").Value = "{'Fuel Tank Sealant': 53.33, 'JOINT COMPOUND': 158.01, 'NUT': 6.13}"
$ws.Range("G58").Value = "{'Fuel Tank Sealant': 43.07, 'PIN': 0.23}"
$ws.Range("G59").Value = 0
$ws.Range("G60").Value = "{'SEAL': 215.28}"
$ws.Range("G61").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("G63").Value = 0
